# RenardESPBOM.xlsx edit: "Changed Capacitor in BOM" / "Bigger!"
#
# - Swap the Capacitor Pol part number (row 3) for a new Digikey part,
#   and update its unit price / price-per-100 figures.
# - Add two new columns (J = "100 QTY price/ea", K = "100 qty price tot")
#   of bulk-pricing data for the existing BOM rows, plus two extra rows
#   (24, 25) of bulk pricing for other items.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Capacitor Pol gets a new part number / pricing -----------------
$ws.Range("B3").Value = "1189-2375-ND"
$ws.Range("E3").Value = 0.29
$ws.Range("G3").Value = 0.11008

# --- New column J/K bulk-pricing data --------------------------------------

# Row 2
$ws.Range("J2").Value = 0.416
$ws.Range("K2").Formula = "=J2*D2"

# Rows 3-11 share one formula pattern (=Jn*Dn) applied as a single range so
# Excel stores it as a shared formula, just like the source workbook.
$ws.Range("J3").Value = 0.156
$ws.Range("J4").Value = 0.202
$ws.Range("J5").Value = 0.595
$ws.Range("J6").Value = 0.744
$ws.Range("J7").Value = 0.749
$ws.Range("J8").Value = 0.053
$ws.Range("J9").Value = 9.61
$ws.Range("J10").Value = 2.34
$ws.Range("J11").Value = 0.052
$ws.Range("K3:K11").Formula = "=J3*D3"

# Row 14 - plain values, no formula
$ws.Range("J14").Value = 1.3
$ws.Range("K14").Value = 1.3

# Row 17 - plain value, no formula, no J17
$ws.Range("K17").Value = 5.56

# Rows 24 and 25 - new rows with their own bulk pricing formula (=Jn/4)
$ws.Range("J24").Value = 14.8
$ws.Range("K24").Formula = "=J24/4"
$ws.Range("J25").Value = 15.3
$ws.Range("K25").Formula = "=J25/4"

# --- Cosmetic: the saved selection moved to G5 ------------------------------
$ws.Range("G5").Select() | Out-Null

$wb.Save()
